# Combine Florenceville and Bristol prior to 2009
# Florenceville, Bristol, and Florenceville-Bristol (pre-2009 combined record)
# rows are removed from the inconsistent_munis table; the remaining rows
# shift up to fill the gap and the table/autoFilter/dimension shrink
# accordingly (handled automatically by Excel when rows are deleted).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows in descending order so row numbers of not-yet-deleted rows
# stay valid while we work.
$ws.Rows("8").Delete()   # Florenceville-Bristol (pre-2009 combined record)
$ws.Rows("7").Delete()   # Florenceville
$ws.Rows("4").Delete()   # Bristol
